$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

$row = 7

$ws.Cells.Item($row, 1).Value = 8

# Column B holds a date-like string ("2026-02-16") that must stay text,
# not get auto-converted into a date serial number, so force text format
# before assigning it (mirrors existing rows which store it as text).
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "2026-02-16"

$ws.Cells.Item($row, 3).Value = "21:21:59"
$ws.Cells.Item($row, 4).Value = "leadlag"
$ws.Cells.Item($row, 5).Value = "DOWN"
$ws.Cells.Item($row, 6).Value = 69364.49000000001
# Column G (Exit Price) stays blank - trade is still OPEN.
$ws.Cells.Item($row, 8).Value = "OPEN"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0.7448
$ws.Cells.Item($row, 12).Value = "Binance leading with -0.074% move"
# Column M (Exit Reason) stays blank - trade is still OPEN.
$ws.Cells.Item($row, 14).Value = 0
